$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 2-4 get some of their values corrected / rows 5 onward are brand-new
# (train circuit "obieg" number -> rolling-stock type table, extended from 4 to 15 entries)
$data = @(
  @(1,  "815 01", "EN71"),
  @(2,  "815 02", "ED72Ac"),
  @(3,  "815 03", "ED72Ac"),
  @(4,  "815 04", "2x EN57FPS"),
  @(5,  "815 05", "EN57FPS"),
  @(6,  "815 06", "EN57FPS"),
  @(7,  "815 07", "EN57FPS"),
  @(8,  "815 08", "EN57ALwKM"),
  @(9,  "815 09", "EN57ALwKM"),
  @(10, "815 10", "EN57ALwKM"),
  @(11, "815 11", "EN57"),
  @(12, "815 12", "EN57"),
  @(13, "815 13", "EN57"),
  @(14, "815 14", "EN57"),
  @(15, "815 15", "EN57")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r++
}

$ws.Range("D25").Select()
